$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header F1
$ws.Range("F1").Value = "Dia diem"

# Update row 2 values with control-name placeholders
$ws.Range("B2").Value = ".!entry"
$ws.Range("C2").Value = ".!combobox"
$ws.Range("D2").Value = ".!entry2"
$ws.Range("E2").Value = ".!checkbutton1"
$ws.Range("F2").Value = ".!entry3"
$ws.Range("G2").Value = ".!spinbox"

# Update row 3: clear B3, C3, D3; set E3; clear F3, G3
$ws.Range("B3").Value = ""
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ".!checkbutton2"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = ""

# Remove rows 4 and 5 entirely
$ws.Rows("4:5").Delete()
